$wb = $excel.ActiveWorkbook

# Rename "batt_portfolio" -> "component_portfolio" and update its view/selection.
$battSheet = $wb.Worksheets.Item("batt_portfolio")
$battSheet.Name = "component_portfolio"
$battSheet.Select()
$battSheet.Range("F41").Select()

# Rename "enr_tec_correspondance" -> "enr_tec_correspondence", make it the
# active/selected tab, and update its selection.
$corrSheet = $wb.Worksheets.Item("enr_tec_correspondance")
$corrSheet.Name = "enr_tec_correspondence"
$corrSheet.Activate()
$corrSheet.Select()
$corrSheet.Range("J47").Select()
